$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Data")
$wsCodebook = $wb.Worksheets.Item("Codebook")

# --- Data sheet: rename the "Eye color" / "Waist" columns to
# "Inseam" (numeric, cm) and "Hair Color" (categorical) ---
$wsData.Range("D1").Value = "Inseam"
$wsData.Range("E1").Value = "Hair Color"

# New Inseam (cm) and Hair Color values for rows 2..15
$inseam = @{
    2  = 81
    3  = 74
    4  = 55
    5  = 91
    6  = 95
    7  = 89
    8  = 68
    9  = 62
    10 = 73
    11 = 49
    12 = 53
    13 = 50
    14 = 52
    15 = 65
}
$hairColor = @{
    2  = "black"
    3  = "blond"
    4  = "brown"
    5  = "l brown"
    6  = "other"
    7  = "black"
    8  = "other"
    9  = "blond"
    10 = "brown"
    11 = "other"
    12 = "other"
    13 = "l brown"
    14 = "d brown"
    15 = "white"
}

foreach ($row in 2..15) {
    $wsData.Range("D$row").Value = $inseam[$row]
    $wsData.Range("E$row").Value = $hairColor[$row]
}

# --- Codebook sheet: remove the "Eye Color" and "Waist" variable
# definitions (rows 5 and 6); everything else is unchanged ---
$wsCodebook.Rows.Item(5).Delete()
$wsCodebook.Rows.Item(5).Delete()

# Leave the Codebook selection resting on the last remaining row
# (A4), then return to the Data sheet as the active tab.
$wsCodebook.Activate()
$wsCodebook.Range("A4").Select()

# --- Selection / active sheet bookkeeping, matching the saved
# state after editing the last Hair Color cell on the Data sheet ---
$wsData.Activate()
$wsData.Range("E15").Select()
